# cambio de plan 4-6-2015
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan de estudio")

# --- Row 3: "Dias para la entrega" 4 -> 5 ---
$ws.Range("N3").Value = 5

# --- Row 5: "Dias para la entrega" 4 -> 5 ---
$ws.Range("N5").Value = 5

# --- Row 6: "Dias para la entrega" 4 -> 5 ---
$ws.Range("N6").Value = 5

# --- Row 7: "Dias para la entrega" 5 -> 10 ---
$ws.Range("N7").Value = 10

# --- Row 8: "Dias para la entrega" 5 -> 10 ---
$ws.Range("N8").Value = 10

# --- Row 9: "Dias para la entrega" 5 -> 10 ---
$ws.Range("N9").Value = 10

# --- Row 10: "Porcentaje"/L10 3 -> 0 ---
$ws.Range("L10").Value = 0

# --- Row 11: fill in the previously-empty "Trabajos pendientes" entry ---
# O11 needs the "Terminado" (done) status look, same as O10 - copy that cell's
# format down before setting the values.
$ws.Range("O10").Copy()
$ws.Range("O11").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("J11").Value = "Emprendedurismo"
$ws.Range("K11").Value = "viernes 5 de agosto"
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = "Mapa de empatia"
$ws.Range("N11").Value = 3
$ws.Range("O11").Value = "Terminado"

# --- move the active selection to O7, matching the saved view state ---
[void]$ws.Range("O7").Select()
